$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data cell T2 (was 252466) to the new value 256389
$ws.Range("T2").Value = 256389

# Move the active selection from T3 to T2
$ws.Range("T2").Select()
